# Fix che, omsk, krd
# - correct the id in A411 (822 -> 821)
# - insert a new row for Soueast S09 at row 412, shifting subsequent rows down

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mis-ordered id for Soueast S07 (row 411) - keep stored as text like the rest of column A
$ws.Range("A411").NumberFormat = "@"
$ws.Range("A411").Value = "821"

# Insert a new row at 412 (pushes old row 412.. down by one)
$ws.Rows.Item(412).EntireRow.Insert()

# Populate the newly inserted row 412 with Soueast S09 data
$ws.Range("A412").NumberFormat = "@"
$ws.Range("A412").Value = "822"
$ws.Range("B412").Value = "Soueast"
$ws.Range("C412").Value = "S09"
$ws.Range("D412").Value = 2050000
$ws.Range("E412").Value = "https://spb.carso.ru/soueast/s09"
$ws.Range("J412").Value = 2050000
$ws.Range("K412").Value = "https://spb.carso.ru/soueast/s09"
